$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "test.user@test.com"
$ws.Range("D2").Value = "Test User"
$ws.Range("F2").Value = "0 B"

$ws.Range("F15").Select()
